$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename two study_id entries to more specific names (Maxwell et al. 2023 hook update)
$ws.Range("A23").Value = "Yando_et_al_2016_marsh"
$ws.Range("A12").Value = "Miller_et_al_2022_Scotland"
